$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.256.79"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "3.206.08"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'608.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "'156.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.205.14"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'0.553"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("D11").Value = "'5.66"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.47%  "
$ws.Range("D12").Value = "'0.504"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("D14").Value = "'38.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "3.728.69"
$ws.Range("E15").Value = "  +0.38%  "
$ws.Range("D16").Value = "66.385.74"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'7.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").Value = "3.207.03"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("E19").Value = "  +1.50%  "
$ws.Range("D20").Value = "'508.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("D22").Value = "'0.735"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").Value = "'8.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("D24").Value = "'14.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'85.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D27").Value = "'3.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'9.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.03%  "
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("E30").Value = "  +41.73%  "
$ws.Range("D31").Value = "'2.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'7.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.64%  "
$ws.Range("D33").Value = "'28.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  -4.79%  "
$ws.Range("D36").Value = "'6.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.57%  "
$ws.Range("D37").Value = "'502.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").Value = "'55.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").Value = "0.0₃0776"
$ws.Range("E39").Value = "  +15.77%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.0422"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.131"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("D42").Value = "'3.07"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").Value = "'8.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.22%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.906.90"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'28.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").Value = "'122.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.12%  "
